$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: Nikita Dane, 1111111111, French
$ws.Range("A2").Value = "Nikita Dane"
$ws.Range("B2").Value = 1111111111
$ws.Range("C2").Value = "French"

# Update row 3: Daphne Fong, 2222222222, Arabic
$ws.Range("A3").Value = "Daphne Fong"
$ws.Range("B3").Value = 2222222222
$ws.Range("C3").Value = "Arabic"

# Set column B width to (best) fit its contents
$ws.Columns.Item(2).ColumnWidth = 10.33

# Update the active selection to C2
$ws.Range("C2").Select()
